$d = $word.ActiveDocument

$d.Content.Find.Execute("431÷4=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "453÷7=64, 5", 2)
$d.Content.Find.Execute("572÷4=143, 0", $true, $false, $false, $false, $false, $true, 1, $false, "289÷8=36, 1", 2)
$d.Content.Find.Execute("437÷8=54, 5", $true, $false, $false, $false, $false, $true, 1, $false, "927÷3=309, 0", 2)
$d.Content.Find.Execute("408÷7=58, 2", $true, $false, $false, $false, $false, $true, 1, $false, "380÷5=76, 0", 2)
$d.Content.Find.Execute("382÷4=95, 2", $true, $false, $false, $false, $false, $true, 1, $false, "653÷2=326, 1", 2)
$d.Content.Find.Execute("957÷9=106, 3", $true, $false, $false, $false, $false, $true, 1, $false, "556÷4=139, 0", 2)
$d.Content.Find.Execute("508÷2=254, 0", $true, $false, $false, $false, $false, $true, 1, $false, "963÷9=107, 0", 2)
$d.Content.Find.Execute("253÷4=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "135÷4=33, 3", 2)
$d.Content.Find.Execute("577÷3=192, 1", $true, $false, $false, $false, $false, $true, 1, $false, "172÷3=57, 1", 2)
$d.Content.Find.Execute("234÷6=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "420÷8=52, 4", 2)
$d.Content.Find.Execute("462÷9=51, 3", $true, $false, $false, $false, $false, $true, 1, $false, "260÷3=86, 2", 2)
$d.Content.Find.Execute("878÷3=292, 2", $true, $false, $false, $false, $false, $true, 1, $false, "341÷9=37, 8", 2)
$d.Content.Find.Execute("154÷4=38, 2", $true, $false, $false, $false, $false, $true, 1, $false, "878÷4=219, 2", 2)
$d.Content.Find.Execute("926÷8=115, 6", $true, $false, $false, $false, $false, $true, 1, $false, "130÷3=43, 1", 2)
$d.Content.Find.Execute("883÷6=147, 1", $true, $false, $false, $false, $false, $true, 1, $false, "940÷9=104, 4", 2)
$d.Content.Find.Execute("514÷2=257, 0", $true, $false, $false, $false, $false, $true, 1, $false, "446÷2=223, 0", 2)
$d.Content.Find.Execute("523÷7=74, 5", $true, $false, $false, $false, $false, $true, 1, $false, "963÷9=107, 0", 2)
$d.Content.Find.Execute("453÷8=56, 5", $true, $false, $false, $false, $false, $true, 1, $false, "685÷7=97, 6", 2)
$d.Content.Find.Execute("454÷7=64, 6", $true, $false, $false, $false, $false, $true, 1, $false, "769÷3=256, 1", 2)
$d.Content.Find.Execute("297÷2=148, 1", $true, $false, $false, $false, $false, $true, 1, $false, "380÷2=190, 0", 2)
$d.Content.Find.Execute("815÷6=135, 5", $true, $false, $false, $false, $false, $true, 1, $false, "899÷8=112, 3", 2)
$d.Content.Find.Execute("755÷3=251, 2", $true, $false, $false, $false, $false, $true, 1, $false, "125÷3=41, 2", 2)
$d.Content.Find.Execute("956÷2=478, 0", $true, $false, $false, $false, $false, $true, 1, $false, "612÷2=306, 0", 2)
$d.Content.Find.Execute("950÷2=475, 0", $true, $false, $false, $false, $false, $true, 1, $false, "494÷8=61, 6", 2)
$d.Content.Find.Execute("770÷9=85, 5", $true, $false, $false, $false, $false, $true, 1, $false, "180÷4=45, 0", 2)
